$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Remove the four discontinued category rows (bottom-to-top so row
#     numbers of not-yet-deleted rows stay stable) ---
$ws.Range("A14:F14").EntireRow.Delete()   # PANELES PVC
$ws.Range("A13:F13").EntireRow.Delete()   # PANELES PU
$ws.Range("A9:F9").EntireRow.Delete()     # LED
$ws.Range("A5:F5").EntireRow.Delete()     # GRANITO

# --- Refresh the remaining category rows with the updated figures ---
$ws.Range("B2").Value = "240X120 PORCELANATO"
$ws.Range("C2").Value = 2892.20588040374
$ws.Range("D2").Value = 388.8
$ws.Range("E2").Value = 2503.40588040374
$ws.Range("F2").Value = 0.1344302639844315

$ws.Range("B3").Value = "240X80 PORCELANATO"
$ws.Range("C3").Value = 20387.4774217135
$ws.Range("D3").Value = 345.99
$ws.Range("E3").Value = 20041.4874217135
$ws.Range("F3").Value = 0.01697071162082595

$ws.Range("B4").Value = "FREGADEROS DE COCINA"
$ws.Range("C4").Value = 782.417163948959
$ws.Range("D4").Value = 65.79000000000001
$ws.Range("E4").Value = 716.627163948959
$ws.Range("F4").Value = 0.08408557868023946

$ws.Range("B5").Value = "GRIFERIAS"
$ws.Range("C5").Value = 150
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 150
$ws.Range("F5").Value = 0

$ws.Range("B6").Value = "INODOROS"
$ws.Range("C6").Value = 2326.06694516821
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 2326.06694516821
$ws.Range("F6").Value = 0

$ws.Range("B7").Value = "LAVABOS"
$ws.Range("C7").Value = 886.711016287574
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 886.711016287574
$ws.Range("F7").Value = 0

$ws.Range("B8").Value = "NO RESURTIBLES"
$ws.Range("C8").Value = 448.80162917203
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 448.80162917203
$ws.Range("F8").Value = 0

$ws.Range("B9").Value = "OTROS"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

$ws.Range("B10").Value = "PANELES DECORATIVOS"
$ws.Range("C10").Value = 2716.75588474074
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 2716.75588474074
$ws.Range("F10").Value = 0

$ws.Range("B11").Value = "PIEDRA SINTERIZADA"
$ws.Range("C11").Value = 19573.0602492497
$ws.Range("D11").Value = 2518.21
$ws.Range("E11").Value = 17054.8502492497
$ws.Range("F11").Value = 0.1286569380532373

$ws.Range("B12").Value = "PORCELANATO"
$ws.Range("C12").Value = 47134.2631579098
$ws.Range("D12").Value = 6524.65
$ws.Range("E12").Value = 40609.6131579098
$ws.Range("F12").Value = 0.1384269014271218

$ws.Range("B13").Value = "PUERTAS DE SEGURIDAD"
$ws.Range("C13").Value = 1110.43665120341
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1110.43665120341
$ws.Range("F13").Value = 0

$ws.Range("B14").Value = "SAL SOLUBLE"
$ws.Range("C14").Value = 489.803925295083
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 489.803925295083
$ws.Range("F14").Value = 0

# --- TOTAL row (moved up from row 19 to row 15) ---
$ws.Range("B15").Value = "TOTAL"
$ws.Range("C15").Value = 98897.99992509274
$ws.Range("D15").Value = 9843.439999999999
$ws.Range("E15").Value = 89054.55992509275
$ws.Range("F15").Value = 0.09953123427628073

# --- Updated column widths ---
# (COM's ColumnWidth stores 5/6 of a character narrower than the saved
#  OOXML <col width> value, so compensate by +5/6 to land on the exact
#  target width once persisted.)
$ws.Columns.Item(3).ColumnWidth = 21.1666666666667
$ws.Columns.Item(4).ColumnWidth = 12.1666666666667
$ws.Columns.Item(5).ColumnWidth = 21.1666666666667
$ws.Columns.Item(6).ColumnWidth = 24.1666666666667
